$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 22:16"

# Insert fresh data for Asturias (now sorted above Malaga by total cases)
# and shift the following rows (Malaga, Ciudad Real, Alacant/Alicante, Toledo) down by one.
$ws.Range("A11").Value = "Asturias"
$ws.Range("B11").Value = 486
$ws.Range("C11").Value = 12
$ws.Range("D11").Value = 469
$ws.Range("E11").Value = 5

$ws.Range("A12").Value = "Malaga"
$ws.Range("B12").Value = 424
$ws.Range("C12").Value = 72
$ws.Range("D12").Value = 406
$ws.Range("E12").Value = 18

$ws.Range("A13").Value = "Ciudad Real"
$ws.Range("B13").Value = 400
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 364
$ws.Range("E13").Value = 28

$ws.Range("A14").Value = "Alacant/Alicante"
$ws.Range("B14").Value = 372
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 348
$ws.Range("E14").Value = 17

$ws.Range("A15").Value = "Toledo"
$ws.Range("B15").Value = 370
$ws.Range("C15").Value = 15
$ws.Range("D15").Value = 336
$ws.Range("E15").Value = 19

# Swap order of Arroyo de la Luz / La Palma (tied totals, alphabetised/reordered)
$ws.Range("A56").Value = "Arroyo de la Luz"
$ws.Range("A57").Value = "La Palma"
